$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows for "Corazón de apio" dated 2022-03-03 (row 4) and
# 2021-06-23 (row 5) were swapped (all columns D..Q except the
# unchanged O/R values), effectively reordering the two weekly
# observations chronologically.

# --- Row 4: set to old row 5 values ---
$ws.Range("D4").Value = 44370
$ws.Range("I4").Value = "Segunda"
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 1000
$ws.Range("L4").Value = 1200
$ws.Range("M4").Value = 1080
$ws.Range("N4").Value = "`$/docena de matas"
$ws.Range("P4").Value = 180
$ws.Range("Q4").Value = 6

# --- Row 5: set to old row 4 values ---
$ws.Range("D5").Value = 44623
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 300
$ws.Range("K5").Value = 1800
$ws.Range("L5").Value = 2000
$ws.Range("M5").Value = 1900
$ws.Range("N5").Value = "`$/paquete"
$ws.Range("P5").Value = 1900
$ws.Range("Q5").Value = 1
